# Apply updated dSF (column F) values for the 2024 lynch_daniel data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7  = 7
    13 = -9
    14 = -8
    15 = 1
    17 = 2
    18 = 1
    20 = 3
    21 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
